# Apply fixes for two hardware issues:
#  1) Buck current measuring negative -> rescale HvBuckBootstrap2 test row (row 22)
#  2) FLT_OUT signal not behaving correctly -> swap Active/Inactive rows for
#     SPM (29/30) and FLTOUT (33/34) test rows
# Also widen column B on the Report sheet, update the active selections, and
# correct a couple of Quantities sheet bounds.

$wb = $excel.ActiveWorkbook

$wsReport = $wb.Worksheets.Item("Report")
$wsQuantities = $wb.Worksheets.Item("Quantities")

# --- Report sheet: column B width ---
# Stored OOXML width = ColumnWidth + 0.8333333333333334, so back the COM
# property off by that offset to land on a stored width of exactly 53.
$wsReport.Columns.Item(2).ColumnWidth = 53 - 0.8333333333333334

# --- Report sheet: row 22 (HvBuckBootstrap2) ---
$wsReport.Range("D22").Value = 235
$wsReport.Range("E22").Value = 250
$wsReport.Range("F22").Value = 240

# --- Report sheet: row 29 (SpmActive) ---
$wsReport.Range("D29").Value = 4.8
$wsReport.Range("E29").Value = 5.0999999999999996
$wsReport.Range("F29").Value = 5

# --- Report sheet: row 30 (SpmInactive) ---
$wsReport.Range("D30").Value = 0
$wsReport.Range("E30").Value = 0.3
$wsReport.Range("F30").Value = 0.2

# --- Report sheet: row 33 (FltOutActive) ---
$wsReport.Range("D33").Value = 4.8
$wsReport.Range("E33").Value = 5.0999999999999996
$wsReport.Range("F33").Value = 5

# --- Report sheet: row 34 (FltOutInactive) ---
$wsReport.Range("D34").Value = 0
$wsReport.Range("E34").Value = 0.3
$wsReport.Range("F34").Value = 0.2

# --- Quantities sheet: row 11 ---
$wsQuantities.Range("C11").Value = 235
$wsQuantities.Range("D11").Value = 245

# --- Update active selections to match author's final cursor positions ---
$wsReport.Activate()
[void]$wsReport.Range("I35").Select()

$wsQuantities.Activate()
[void]$wsQuantities.Range("E20").Select()

$wsReport.Activate()
